# Overall_Rebate_Efficiency.xlsx edit:
#  - Update the PSA_LOLO sheet's data row (A2, B2) with new values
#  - Switch the active/selected sheet from OverallRebateEfficiency to PSA_LOLO
#  - Leave behind a selection of B23 on OverallRebateEfficiency (no longer the
#    active sheet) and keep A3 selected on PSA_LOLO

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("OverallRebateEfficiency")
$ws2 = $wb.Worksheets.Item("PSA_LOLO")

# Last selection left on the first sheet before moving away from it.
$ws1.Range("B23").Select()

# Update the two data values on PSA_LOLO.
$ws2.Range("A2").Value = 43297
$ws2.Range("B2").Value = 25479

# Make PSA_LOLO the active sheet/tab, with A3 selected.
$ws2.Activate()
$ws2.Range("A3").Select()
